$d = $word.ActiveDocument
$rng = $d.Content
$replacements = @(
    @("2024-07-22 Monday", "2024-07-23 Tuesday"),
    @("6+64=70", "29+13=42"),
    @("24+48=72", "19-1=18"),
    @("50+12=62", "21+40=61"),
    @("98-97=1", "79-75=4"),
    @("74-38=36", "92-30=62"),
    @("27-9=18", "81-20=61"),
    @("28+27=55", "22+11=33"),
    @("84-75=9", "25+43=68"),
    @("74-8=66", "89-0=89"),
    @("75-8=67", "16+41=57"),
    @("55+22=77", "15+1=16"),
    @("99-31=68", "37-31=6"),
    @("40-12=28", "69-17=52"),
    @("98-56=42", "11+38=49"),
    @("95-54=41", "2+69=71"),
    @("7+43=50", "76+2=78"),
    @("87-13=74", "13+41=54"),
    @("10+3=13", "89-33=56"),
    @("45-0=45", "69-64=5"),
    @("82-45=37", "27-8=19"),
    @("53+23=76", "91-11=80"),
    @("24-7=17", "81-44=37"),
    @("0+39=39", "75-61=14"),
    @("16+49=65", "15+53=68"),
    @("88-56=32", "32-29=3"),
    @("86-2=84", "22-12=10"),
    @("27+5=32", "40-39=1"),
    @("26+30=56", "26+3=29"),
    @("10+39=49", "85-66=19"),
    @("93+3=96", "30+62=92"),
    @("93-0=93", "19+11=30"),
    @("72-68=4", "86+4=90"),
    @("95-69=26", "44-22=22"),
    @("8+39=47", "6+27=33"),
    @("57+18=75", "14+11=25"),
    @("17+26=43", "39-4=35"),
    @("78-6=72", "31+61=92"),
    @("82-55=27", "37+37=74"),
    @("50-15=35", "82-10=72"),
    @("37+5=42", "56+32=88"),
    @("80-70=10", "77-21=56"),
    @("89-52=37", "39+21=60"),
    @("23+51=74", "24+13=37"),
    @("43+23=66", "15+84=99"),
    @("27-16=11", "75+24=99"),
    @("39+15=54", "38+25=63"),
    @("12+12=24", "16+27=43"),
    @("72-1=71", "61-11=50"),
    @("69+26=95", "94-53=41"),
    @("39+59=98", "27+65=92"),
    @("69+28=97", "77-28=49"),
    @("27+71=98", "34-12=22"),
    @("33+65=98", "77-25=52"),
    @("8-6=2", "29+18=47"),
    @("26-8=18", "60-9=51"),
    @("0+37=37", "4+52=56"),
    @("56+17=73", "9+49=58"),
    @("57+21=78", "4+90=94"),
    @("13+73=86", "33+21=54"),
    @("42+52=94", "67-61=6"),
    @("13-7=6", "45+27=72"),
    @("76+23=99", "48+22=70"),
    @("0+45=45", "88+9=97"),
    @("20+60=80", "83-73=10"),
    @("40-18=22", "19+55=74"),
    @("23+38=61", "1+78=79"),
    @("18+5=23", "89-8=81"),
    @("78-45=33", "97-58=39"),
    @("68-61=7", "98-90=8"),
    @("51-48=3", "18+43=61"),
    @("8-6=2", "84-20=64"),
    @("19+10=29", "99-87=12"),
    @("94-45=49", "46-43=3"),
    @("45+34=79", "58+22=80"),
    @("43-8=35", "5+60=65"),
    @("80+14=94", "82-57=25"),
    @("16-13=3", "27+63=90"),
    @("51+30=81", "58-50=8"),
    @("21-8=13", "46-15=31"),
    @("66+14=80", "66-61=5"),
    @("59+9=68", "59-19=40"),
    @("67-12=55", "91-45=46"),
    @("52+9=61", "70-68=2"),
    @("42+34=76", "92-89=3"),
    @("75-16=59", "72-71=1"),
    @("24-18=6", "32+52=84"),
    @("99-24=75", "23+22=45"),
    @("24+56=80", "33+49=82"),
    @("23+8=31", "22+63=85"),
    @("84-21=63", "1+22=23"),
    @("15+33=48", "98-1=97"),
    @("5+37=42", "55-29=26"),
    @("94-85=9", "22+75=97"),
    @("37-7=30", "29+17=46"),
    @("53+34=87", "14+74=88"),
    @("74-39=35", "99-87=12"),
    @("16+52=68", "78-52=26"),
    @("72+20=92", "29-27=2"),
    @("16+8=24", "12+84=96"),
    @("44+37=81", "56-14=42"),
)

$failCount = 0
foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $ok = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 1)
    if (-not $ok) {
        $failCount = $failCount + 1
        Write-Host "FAILED to replace: $old -> $new"
    }
}
Write-Host "Done. Failures: $failCount"
